# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45202 (2023-10-03) to 45203 (2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
